# Apply the batch_test iteration-12 data refresh: quantities (column C) and
# their corresponding upto-date amounts (column G) change for rows 8-17,
# which cascade into the Grand Total / Net Payable rows (19 and 21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Force the cell to keep a literal text value (matching the workbook's
    # existing convention of storing formatted amounts like "24576.00" as
    # text) instead of letting Excel auto-coerce the numeric-looking string
    # into a real number. Revert the cell style back to Normal afterwards so
    # no stray formatting is left behind on the cell.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

# Row 8 - Qty executed upto date
$ws.Range("C8").Value = 99

# Row 9 - Short point (up to 3 mtr.)
$ws.Range("C9").Value = 96
Set-TextValue "G9" "24576.00"

# Row 10 - Medium point (up to 6 mtr.)
$ws.Range("C10").Value = 92
Set-TextValue "G10" "43424.00"

# Row 11 - Long point (up to 10 mtr.)
$ws.Range("C11").Value = 88
Set-TextValue "G11" "58256.00"

# Row 12 - Qty executed upto date
$ws.Range("C12").Value = 99

# Row 13 - On board
$ws.Range("C13").Value = 96
Set-TextValue "G13" "13056.00"

# Row 14 - P & F 6 amp switch
$ws.Range("C14").Value = 39
Set-TextValue "G14" "897.00"

# Row 15 - Total
$ws.Range("C15").Value = 10

# Row 16 - Add Tender Premium
$ws.Range("C16").Value = 53

# Row 17 - Grand Total
$ws.Range("C17").Value = 12

# Row 19 - Grand Total Rs.
Set-TextValue "G19" "140209.00"
Set-TextValue "H19" "140209.00"

# Row 21 - NET PAYABLE AMOUNT Rs.
Set-TextValue "G21" "140209.00"
Set-TextValue "H21" "140209.00"
